$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week dates) ---
$ws.Range("A8").Characters(21, 2).Text = "26"
$ws.Range("C9").Characters(27, 9).Text = "6/23/2025"
$ws.Range("C9").Characters(47, 9).Text = "6/29/2025"

# --- Numeric data cell updates (crime statistics table) ---
$ws.Range("C15").Value = 1
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 200
$ws.Range("I15").Value = 18
$ws.Range("K15").Value = 125
$ws.Range("L15").Value = 260
$ws.Range("D16").Value = 1
$ws.Range("F16").Value = 7
$ws.Range("H16").Value = 16.666666666666
$ws.Range("J16").Value = 42
$ws.Range("K16").Value = -9.523809523809
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 17
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 109
$ws.Range("J17").Value = 106
$ws.Range("K17").Value = 2.830188679245
$ws.Range("L17").Value = -9.166666666666
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -45.454545454545
$ws.Range("I18").Value = 39
$ws.Range("J18").Value = 39
$ws.Range("K18").Value = 0
$ws.Range("F19").Value = 29
$ws.Range("G19").Value = 21
$ws.Range("H19").Value = 38.095238095238
$ws.Range("I19").Value = 194
$ws.Range("J19").Value = 170
$ws.Range("K19").Value = 14.117647058823
$ws.Range("L19").Value = -3.960396039603
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -75
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = -66.666666666666
$ws.Range("I20").Value = 29
$ws.Range("J20").Value = 41
$ws.Range("K20").Value = -29.268292682926
$ws.Range("L20").Value = -46.296296296296
$ws.Range("C21").Value = 13
$ws.Range("D21").Value = 10
$ws.Range("E21").Value = 30
$ws.Range("F21").Value = 66
$ws.Range("G21").Value = 68
$ws.Range("H21").Value = -2.941176470588
$ws.Range("I21").Value = 427
$ws.Range("J21").Value = 406
$ws.Range("K21").Value = 5.172413793103
$ws.Range("L21").Value = -4.044943820224
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 25
$ws.Range("E24").Value = -8
$ws.Range("F24").Value = 97
$ws.Range("G24").Value = 109
$ws.Range("H24").Value = -11.009174311926
$ws.Range("I24").Value = 723
$ws.Range("J24").Value = 732
$ws.Range("K24").Value = -1.229508196721
$ws.Range("L24").Value = 7.589285714285
$ws.Range("C25").Value = 18
$ws.Range("D25").Value = 17
$ws.Range("E25").Value = 5.882352941176
$ws.Range("F25").Value = 69
$ws.Range("G25").Value = 73
$ws.Range("H25").Value = -5.479452054794
$ws.Range("I25").Value = 501
$ws.Range("J25").Value = 505
$ws.Range("K25").Value = -0.792079207920
$ws.Range("L25").Value = 24.009900990099
$ws.Range("C26").Value = 11
$ws.Range("E26").Value = -8.333333333333
$ws.Range("F26").Value = 50
$ws.Range("G26").Value = 46
$ws.Range("H26").Value = 8.695652173913
$ws.Range("I26").Value = 288
$ws.Range("J26").Value = 297
$ws.Range("K26").Value = -3.030303030303
$ws.Range("L26").Value = 5.882352941176
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 22
$ws.Range("K27").Value = 37.5
$ws.Range("L27").Value = 144.444444444444
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = -16.666666666666
$ws.Range("I28").Value = 31
$ws.Range("J28").Value = 30
$ws.Range("K28").Value = 3.333333333333
$ws.Range("L28").Value = -13.888888888888

# --- Cells converted from numeric 0 / blank-style placeholders to text placeholders ("0" / "***.*") ---
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C23").PasteSpecial(-4122) | Out-Null
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0"
$ws.Range("D14").Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4122) | Out-Null
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "***.*"
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E23").PasteSpecial(-4122) | Out-Null
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0"
$ws.Range("D14").Copy() | Out-Null
$ws.Range("D33").PasteSpecial(-4122) | Out-Null
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "***.*"
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E33").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
Write-Output "edit complete"
